# Update the JOUR course-requirements sheet to split the single
# "Prerequisites" column into four columns: Prerequisites (C, existing),
# Corequisites (D, new), Concurrent (E, new), Recommended (F, new),
# and shift "Terms Typically Offered" from D to G.
#
# For two courses (JOUR 346 and JOUR 385) the previously-combined
# "Recommended: ..." clause embedded in the Prerequisites text is split
# out into the new Recommended column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header row for the inserted / shifted columns.
$ws.Cells.Item(1, 4).Value = "Corequisites"
$ws.Cells.Item(1, 5).Value = "Concurrent"
$ws.Cells.Item(1, 6).Value = "Recommended"
$ws.Cells.Item(1, 7).Value = "Terms Typically Offered"

$rows = @(
    @{Row=2; C='NA'; D='NA'; E='NA'; F='NA'; G='F'},
    @{Row=3; C='NA'; D='NA'; E='NA'; F='NA'; G='SP'},
    @{Row=4; C='NA'; D='NA'; E='NA'; F='NA'; G='F, W, SP'},
    @{Row=5; C='NA'; D='NA'; E='NA'; F='NA'; G='F, W, SP'},
    @{Row=6; C='NA'; D='NA'; E='NA'; F='NA'; G='W, SP'},
    @{Row=7; C='NA'; D='NA'; E='NA'; F='NA'; G='F, W, SP'},
    @{Row=8; C='Acceptance as a KCPR staff member and consent of instructor.'; D='NA'; E='NA'; F='NA'; G='W'},
    @{Row=9; C='Open to undergraduate students and consent of instructor.'; D='NA'; E='NA'; F='NA'; G='F, W, SP'},
    @{Row=10; C='JOUR 203.'; D='NA'; E='NA'; F='NA'; G='F, W, SP'},
    @{Row=11; C='JOUR 203.'; D='NA'; E='NA'; F='NA'; G='F, W, SP'},
    @{Row=12; C='JOUR 285.'; D='NA'; E='NA'; F='NA'; G='F, W'},
    @{Row=13; C='JOUR 203.'; D='NA'; E='NA'; F='NA'; G='F'},
    @{Row=14; C='Sophomore standing.'; D='NA'; E='NA'; F='NA'; G='F, W'},
    @{Row=15; C='JOUR 220 and consent of instructor.'; D='NA'; E='NA'; F='NA'; G='F, W, SP'},
    @{Row=16; C='Completion of GE Area A3 with a grade of C- or better.'; D='NA'; E='NA'; F='NA'; G='F, W, SP'},
    @{Row=17; C='Completion of GE Area A3 with a grade of C- or better.'; D='NA'; E='NA'; F='NA'; G='F, W, SP'},
    @{Row=18; C='JOUR 203.'; D='NA'; E='NA'; F='NA'; G='F, SP'},
    @{Row=19; C='JOUR 312.'; D='NA'; E='NA'; F='NA'; G='F, W, SP'},
    @{Row=20; C='JOUR 203.'; D='NA'; E='NA'; F='JOUR 285.'; G='W '},
    @{Row=21; C='JOUR 203 and JOUR 333.'; D='NA'; E='NA'; F='NA'; G='SP'},
    @{Row=22; C='JOUR 285; and one of the STAT 130, STAT 217, STAT 218, or STAT 251.'; D='NA'; E='NA'; F='NA'; G='SP'},
    @{Row=23; C='JOUR 304 or JOUR 334.'; D='NA'; E='NA'; F='NA'; G='F, W, SP'},
    @{Row=24; C='JOUR 333, JOUR 346, or JOUR 348.'; D='NA'; E='NA'; F='NA'; G='F, W, SP'},
    @{Row=25; C='JOUR 285 and JOUR 303 or JOUR 346.'; D='NA'; E='NA'; F='NA'; G='SP'},
    @{Row=26; C='Sophomore standing.'; D='NA'; E='NA'; F='BUS 310 or COMS/JOUR 218.'; G='W '},
    @{Row=27; C='JOUR 203.'; D='NA'; E='NA'; F='NA'; G='F, W, SP'},
    @{Row=28; C='Consent of department chair.'; D='NA'; E='NA'; F='NA'; G='F,W,SP,SU'},
    @{Row=29; C='JOUR 203 and junior standing.'; D='NA'; E='NA'; F='NA'; G='SP'},
    @{Row=30; C='JOUR 203.'; D='NA'; E='NA'; F='NA'; G='F, W, SP'},
    @{Row=31; C='JOUR 303 or JOUR 348.'; D='NA'; E='NA'; F='NA'; G='F, W'},
    @{Row=32; C='JOUR 203.'; D='NA'; E='NA'; F='NA'; G='W'},
    @{Row=33; C='JOUR 285.'; D='NA'; E='NA'; F='NA'; G='F, W, SP'},
    @{Row=34; C='Junior standing.'; D='NA'; E='NA'; F='NA'; G='F'},
    @{Row=35; C='JOUR 312 and JOUR 342.'; D='NA'; E='NA'; F='NA'; G='F, W, SP'},
    @{Row=36; C='JOUR 413; or BUS 453, JOUR 331, JOUR 342, and one of the GRC 338, GRC 377, or JOUR 390.'; D='NA'; E='NA'; F='NA'; G='F, W, SP'},
    @{Row=37; C='Junior standing; and completion of GE Area A with a grade of C- or better.'; D='NA'; E='NA'; F='NA'; G='TBD'},
    @{Row=38; C='JOUR 352, JOUR 353, or JOUR 413; and internship coordinator approval.'; D='NA'; E='NA'; F='NA'; G='F,W,SP,SU'},
    @{Row=39; C='Consent of instructor.'; D='NA'; E='NA'; F='NA'; G='F, W, SP'},
    @{Row=40; C='Consent of instructor; Journalism majors only.'; D='NA'; E='NA'; F='NA'; G='F, W, SP'},
    @{Row=41; C='Senior standing.'; D='NA'; E='NA'; F='NA'; G='F, W, SP'},
    @{Row=42; C='Consent of instructor; junior standing.'; D='NA'; E='NA'; F='NA'; G='TBD'},
    @{Row=43; C='Consent of instructor.'; D='NA'; E='NA'; F='NA'; G='TBD'}
)

foreach ($r in $rows) {
    # Column C (Prerequisites) is rewritten for every row: non-breaking
    # spaces between course prefixes and numbers become regular spaces,
    # and embedded "Recommended: ..." / "the following:" clauses are
    # removed now that they live in their own column.
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    # Column D (Corequisites) - new, always "NA" in this data set.
    $ws.Cells.Item($r.Row, 4).Value = $r.D
    # Column E (Concurrent) - new, always "NA" in this data set.
    $ws.Cells.Item($r.Row, 5).Value = $r.E
    # Column F (Recommended) - new; "NA" except where a recommendation
    # was split out of the old Prerequisites text.
    $ws.Cells.Item($r.Row, 6).Value = $r.F
    # Column G (Terms Typically Offered) - shifted from the old column D.
    $ws.Cells.Item($r.Row, 7).Value = $r.G
}
